$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: I2 reordered text
$ws.Range("I2").Value = "{'OR ART. 185 Abs. 1', 'OR ART. 184 Abs. 1', 'OR ART. 214 Abs. 1'}"

# Row 3: I3 now includes 'OR ART. 109 Abs. 1' and reordered
$ws.Range("I3").Value = "{'OR ART. 109 Abs. 1', 'OR ART. 119 Abs. 1', 'OR ART. 185 Abs. 1', 'OR ART. 185 Abs. 2', 'OR ART. 119 Abs. 3', 'OR ART. 119 Abs. 2'}"

# Row 4: I4 filled, K4/L4 set to 0
$ws.Range("I4").Value = "{'OR ART. 259e', 'OR ART. 257g', 'OR ART. 259d', 'OR ART. 266g'}"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0

# Row 5: I5 filled, K5/L5 set to 0
$ws.Range("I5").Value = "{'StGB Art. 179quater', 'OR Art. 19 Abs. 1', 'UrhG Art. 11 Abs. 2', 'ZGB Art. 28a', 'UrhG Art. 13', 'OR Art. 67 Abs. 1', 'OR Art. 41', 'UrhG Art. 20', 'UrhG Art. 28', 'UrhG Art. 15', 'OR Art. 97 Abs. 1', 'OR Art. 423 Abs. 1', 'UrhG Art. 29', 'ZGB Art. 328b', 'UrhG Art. 19 Abs. 1', 'OR Art. 62 Abs. 1', 'UrhG Art. 40', 'ZGB Art. 28 Abs. 1', 'UrhG Art. 36', 'ZGB Art. 28 Abs. 2'}"
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0

# Row 6: I6 filled, K6/L6 set to 0
$ws.Range("I6").Value = "{'OR ART. 18', 'OR ART. 24', 'OR ART. 23', 'OR ART. 97', 'OR ART. 367', 'OR ART. 404', 'OR ART. 364', 'OR ART. 107', 'OR ART. 1', 'OR ART. 363', 'OR ART. 109'}"
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0

# Row 7: I7 filled, K7/L7 set to 0
$ws.Range("I7").Value = "{'OR ART. 330c', 'OR ART. 330b', 'OR ART. 330a', 'OR ART. 335'}"
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

# Row 8: I8 filled, K8/L8 set to 0
$ws.Range("I8").Value = "{'OR ART. 185 Abs. 1', 'OR ART. 184 Abs. 1', 'OR ART. 214 Abs. 1'}"
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0

# Row 9: I9 filled, K9/L9 set to 0
$ws.Range("I9").Value = "{'OR ART. 109 Abs. 1', 'OR ART. 119 Abs. 1', 'OR ART. 185 Abs. 1', 'OR ART. 185 Abs. 2', 'OR ART. 119 Abs. 3', 'OR ART. 119 Abs. 2'}"
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0

# Row 10: I10 filled, K10/L10 set to 0
$ws.Range("I10").Value = "{'OR ART. 185 Abs. 1', 'OR ART. 119 Abs. 1', 'OR ART. 109', 'OR ART. 119 Abs. 2'}"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0

# Row 11: I11 filled, K11/L11 set to 0
$ws.Range("I11").Value = "{'OR ART. 185 Abs. 1', 'OR ART. 185 Abs. 3', 'OR ART. 185 Abs. 2', 'OR ART. 99 Abs. 3'}"
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
